$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3592
$ws.Range("C3").Value = 3804
$ws.Range("C4").Value = 4036
$ws.Range("C5").Value = 4161
$ws.Range("C6").Value = 4669
$ws.Range("C7").Value = 4850
$ws.Range("C8").Value = 4898
$ws.Range("C9").Value = 5005
$ws.Range("C10").Value = 5005
$ws.Range("C11").Value = 5113
$ws.Range("C12").Value = 5113
$ws.Range("C13").Value = 5113
$ws.Range("C14").Value = 5127
$ws.Range("C15").Value = 5163
